$d = $word.ActiveDocument

function Find-Heading1($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

function Insert-BlankBeforeAndPromote($headingText) {
    $p = Find-Heading1 $headingText
    $idx = $p.Index
    $p.Range.InsertParagraphBefore()
    $blank = $d.Paragraphs.Item($idx)
    $heading = $d.Paragraphs.Item($idx + 1)
    $blank.Style = "Normal"
    $blank.Format.SpaceAfter = 0
    $heading.Style = "Heading 2"
}

# ---------------------------------------------------------------------
# 1. "Choosing a database type" : Heading 1 -> Heading 2 (no blank insert)
# ---------------------------------------------------------------------
$p = Find-Heading1 "Choosing a database type"
$p.Style = "Heading 2"

# ---------------------------------------------------------------------
# 2. "Who is going to use it" : insert blank paragraph before + promote
# ---------------------------------------------------------------------
Insert-BlankBeforeAndPromote "Who is going to use it"

# ---------------------------------------------------------------------
# 3. Restructure run split around "insert, update or select"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " as well as whether the user is allowed to insert, update or select. Unfortunately this is insufficient when data from multiple individuals exist in the same table and it is therefore necessary to ad",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# ---------------------------------------------------------------------
# 4. "How is it going to be used" : insert blank paragraph before + promote
# ---------------------------------------------------------------------
Insert-BlankBeforeAndPromote "How is it going to be used"

# ---------------------------------------------------------------------
# 5. "Relational design" : Heading 1 -> Heading 2 (no blank insert)
# ---------------------------------------------------------------------
$p = Find-Heading1 "Relational design"
$p.Style = "Heading 2"
